$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 178: date 2025-11-28 (serial 45989), station "四方坪站充电量(kw)"
$ws.Cells.Item(178, 1).Value = 45989
$ws.Cells.Item(178, 2).Value = "四方坪站充电量(kw)"
$row178 = @(521.91300000000001,1162.5980000000002,425.46899999999999,427.27000000000004,407.31299999999993,390.86799999999999,508.27199999999993,158.73699999999999,220.39,117.94,126.57900000000002,217.36699999999999,705.10900000000004,1369.6659999999999,675.02300000000002,425.589,387.524,259.51499999999999,88.36,207.84400000000002,91.4,58.85,143.24099999999999,54.069999999999993)
for ($i = 0; $i -lt $row178.Length; $i++) {
    $ws.Cells.Item(178, 3 + $i).Value = $row178[$i]
}

# Row 179: date 2025-11-28 (serial 45989), station "高岭站充电量(kw)"
$ws.Cells.Item(179, 1).Value = 45989
$ws.Cells.Item(179, 2).Value = "高岭站充电量(kw)"
$row179 = @(617.5870000000001,662.88900000000001,288.70600000000002,30.835999999999999,11.356,178.328,93.356999999999999,278.29400000000004,178.26000000000002,201.82399999999998,284.94,206.035,369.471,639.0680000000001,191.667,374.971,181.31299999999999,50.766000000000005,28.875,101.553,152.96199999999999,23.867000000000001,63.737000000000002,67.825999999999993)
for ($i = 0; $i -lt $row179.Length; $i++) {
    $ws.Cells.Item(179, 3 + $i).Value = $row179[$i]
}

# Update the active-cell selection to D183 as in the target
$ws.Range("D183").Select() | Out-Null
